$d = $word.ActiveDocument

function Merge-ParagraphText($paraIndex, $newText) {
    # Replaces a paragraph's content with a single run of $newText, merging
    # away any interior proofErr / multi-run splits, while keeping the
    # paragraph's pPr/rPr formatting intact. Works when the paragraph's
    # first child is a <w:r> (i.e. no proofErr sitting *before* the first run).
    $p = $d.Paragraphs.Item($paraIndex)
    $find = $p.Range.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($p.Range.Text.TrimEnd([char]13, [char]7), $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

function Rebuild-ParagraphText($paraIndex, $newText) {
    # Rebuilds a paragraph from scratch by inserting a brand-new paragraph
    # after it with the desired text (inheriting pPr/list formatting via
    # InsertParagraphAfter) and then deleting the original paragraph. This
    # drops any stray proofErr markers that sit *before* the first run,
    # which Find/Replace can't reach.
    $p = $d.Paragraphs.Item($paraIndex)
    $p.Range.InsertParagraphAfter()
    $newp = $d.Paragraphs.Item($paraIndex + 1)
    $newp.Range.Text = $newText
    $p.Range.Delete()
}

# ---------------------------------------------------------------------
# 1) "Filter tags( don't show posts tagged with certain tags)" - merge
#    the 3 runs (split by gramStart/gramEnd proofErr) into a single run.
# ---------------------------------------------------------------------
$apos = [char]8217
$filterTagsText = "Filter tags( don" + $apos + "t show posts tagged with certain tags)"
Merge-ParagraphText 35 $filterTagsText

# ---------------------------------------------------------------------
# 2) Move <w:lastRenderedPageBreak/> from the "...related videos" run to
#    the "...related images" run (paragraphs 39 and 40).
# ---------------------------------------------------------------------
$imagesText = "Get search query related images"
$videosText = "Get search query related videos"

$p39 = $d.Paragraphs.Item(39)
$insertPoint = $d.Range($p39.Range.Start, $p39.Range.Start)
$breakXml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:lastRenderedPageBreak/><w:t>' + $imagesText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($breakXml)

# the insert duplicated the paragraph's text (new run + old run) - drop the
# now-redundant original tail run.
$p39 = $d.Paragraphs.Item(39)
$tailStart = $p39.Range.Start + $imagesText.Length
$tailEnd = $p39.Range.End - 1
$d.Range($tailStart, $tailEnd).Delete()

# Paragraph 40 still carries the old lastRenderedPageBreak - force a clean
# run rewrite (Find/Replace regenerates the run) which drops it while
# keeping the run's rPr (lang) intact.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($videosText, $false, $false, $false, $false, $false, $true, 1, $false, $videosText, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) "Who sent me an ask?" - merge the 3 runs (spellStart/spellEnd around
#    "an") into a single run.
# ---------------------------------------------------------------------
Merge-ParagraphText 49 "Who sent me an ask?"

# ---------------------------------------------------------------------
# 4) "Slugify function" - has a *leading* proofErr (spellStart) before the
#    first run, so Find/Replace can't clear it -> rebuild the paragraph.
# ---------------------------------------------------------------------
Rebuild-ParagraphText 59 "Slugify function"

# ---------------------------------------------------------------------
# 5) "Reblogs -> list of people who reblogged this post" - also has a
#    leading proofErr (spellStart) -> rebuild the paragraph.
# ---------------------------------------------------------------------
$reblogsText = "Reblogs -> list of people who reblogged this post"
Rebuild-ParagraphText 68 $reblogsText

# ---------------------------------------------------------------------
# 6) Insert a new "Create new account" bullet (ListParagraph / numId 4)
#    right after "Filter content", before the "Post" Heading2. Done last
#    so it doesn't shift the paragraph indices used above.
# ---------------------------------------------------------------------
$p12 = $d.Paragraphs.Item(12)
$p12.Range.InsertParagraphAfter()
$newp = $d.Paragraphs.Item(13)
$newp.Range.Text = "Create new account"
